# "add sum of pure money"
# For every stock sheet, add a B6 cell that sums the row of daily net-money
# values (D6:MI6 as a generously-wide open range, matching the existing
# pattern already used for the B8 "average price" sum in each sheet).
# Also updates each sheet's saved selection, and moves the active-tab
# marker from sheet 16 (大金重工) to sheet 18 (普邦股份).

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # New "pure money" total for row 6.
    $ws.Range("B6").Formula = "=SUM(D6:MI6)"

    # Sheets 1 and 2 (达华智能, 中远海发) land their selection on C6;
    # every other sheet lands on B7.
    if ($i -eq 1 -or $i -eq 2) {
        $null = $ws.Range("C6").Select()
    } else {
        $null = $ws.Range("B7").Select()
    }
}

# Move the active tab from sheet 16 to sheet 18 (last one activated wins).
$null = $wb.Worksheets.Item(18).Activate()
$null = $wb.Worksheets.Item(18).Range("B7").Select()
